# Add a new "Spain" market test-data sheet, cloned from the existing
# "Italy" sheet (same layout/styles/merges), then:
#   - rename the clone to "Spain"
#   - update the two market-specific cells (market name + ticket id)
#   - tighten the column widths / row heights to match the new content
#   - move the "active sheet" / selection state from Italy to Spain

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Clone "Italy" right after itself - this carries over all formatting,
# column widths, merged cells, page setup, etc.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Market-specific content
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2038"

# Column widths (characters) - narrower than Italy's
$spain.Columns.Item(1).ColumnWidth = 24.333333333333332
$spain.Columns.Item(2).ColumnWidth = 14.333333333333334
$spain.Columns.Item(4).ColumnWidth = 7.666666666666667

# Rows 3-5 grow to a two-line height once column D narrows (wrapped text)
$spain.Range("A3:D5").RowHeight = 28.8

# Move the selection off Italy (onto its full data range) and make
# the new Spain sheet the active tab with its own selection.
$italy.Range("A1:D10").Select()
$spain.Activate()
$spain.Range("H11").Select()
